$d = $word.ActiveDocument

$replacements = @(
    @("773×7=5411", "883×5=4415"),
    @("776×2=1552", "195×8=1560"),
    @("916×4=3664", "883×3=2649"),
    @("738×9=6642", "435×9=3915"),
    @("762×4=3048", "525×2=1050"),
    @("964×2=1928", "961×5=4805"),
    @("432×8=3456", "767×3=2301"),
    @("184×9=1656", "872×8=6976"),
    @("504×8=4032", "729×8=5832"),
    @("624×5=3120", "251×4=1004"),
    @("695×4=2780", "519×4=2076"),
    @("370×7=2590", "547×7=3829"),
    @("314×9=2826", "811×4=3244"),
    @("780×8=6240", "516×6=3096"),
    @("406×4=1624", "517×4=2068"),
    @("367×7=2569", "519×6=3114"),
    @("278×8=2224", "563×4=2252"),
    @("147×6=882", "621×9=5589"),
    @("619×9=5571", "101×8=808"),
    @("563×3=1689", "224×9=2016"),
    @("362×9=3258", "722×4=2888"),
    @("814×3=2442", "228×3=684"),
    @("603×2=1206", "719×5=3595"),
    @("396×2=792", "122×3=366"),
    @("631×3=1893", "715×3=2145")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
